$wb = $excel.ActiveWorkbook

# --- "info" sheet: update submitter id / names on row 1 ---
$infoSheet = $wb.Worksheets.Item("info")
# Format as text first so the leading zeros in "000001" survive (otherwise
# Excel would coerce the numeric-looking string into the number 1).
$infoSheet.Range("A1").NumberFormat = "@"
$infoSheet.Range("A1").Value = "000001"
$infoSheet.Range("B1").Value = "hong"
$infoSheet.Range("C1").Value = "gil"

# --- "items" sheet: update item names/prices (rows 2-5) ---
$itemsSheet = $wb.Worksheets.Item("items")

$itemsSheet.Range("A2").Value = "노트"
$itemsSheet.Range("C2").Value = 600
$itemsSheet.Range("E2").Value = 600

$itemsSheet.Range("A3").Value = "육개장(컵)"
$itemsSheet.Range("C3").Value = 1050
$itemsSheet.Range("E3").Value = 1050

$itemsSheet.Range("A4").Value = "접시100"
$itemsSheet.Range("C4").Value = 1000
$itemsSheet.Range("E4").Value = 1000

$itemsSheet.Range("A5").Value = "육개장(컵)"
$itemsSheet.Range("C5").Value = 1050
$itemsSheet.Range("E5").Value = 1050
